$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(186).Insert()

$ws.Cells.Item(186, 1).Value = 10
$ws.Cells.Item(186, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(186, 3).Value = "La Araucanía"
$ws.Cells.Item(186, 4).Value = 44680
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 6).Value = 100112044
$ws.Cells.Item(186, 7).Value = "Perejil"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 30
$ws.Cells.Item(186, 11).Value = 4000
$ws.Cells.Item(186, 12).Value = 4000
$ws.Cells.Item(186, 13).Value = 4000
$ws.Cells.Item(186, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(186, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(186, 16).Value = 1333
$ws.Cells.Item(186, 17).Value = 3
$ws.Cells.Item(186, 18).Value = "Hortaliza"
